$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 64, shifting old rows 64:67 down to 66:69
$ws.Range("A64:A65").EntireRow.Insert()

# Fill new row 64 (Banquete, week of 2021-11-16)
$ws.Range("A64").Value = 12
$ws.Range("B64").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C64").Value = "Metropolitana"
$ws.Range("D64").Value = 44516
$ws.Range("D64").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E64").Value = 13
$ws.Range("F64").Value = 300000000
$ws.Range("G64").Value = "Espárragos"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Banquete"
$ws.Range("J64").Value = 410
$ws.Range("K64").Value = 1200
$ws.Range("L64").Value = 1200
$ws.Range("M64").Value = 1200
$ws.Range("N64").Value = "$/kilo"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 1200
$ws.Range("Q64").Value = 1
$ws.Range("R64").Value = "Hortaliza"

# Fill new row 65 (Primera, week of 2021-11-16)
$ws.Range("A65").Value = 12
$ws.Range("B65").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C65").Value = "Metropolitana"
$ws.Range("D65").Value = 44516
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E65").Value = 13
$ws.Range("F65").Value = 300000000
$ws.Range("G65").Value = "Espárragos"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 350
$ws.Range("K65").Value = 1000
$ws.Range("L65").Value = 1000
$ws.Range("M65").Value = 1000
$ws.Range("N65").Value = "$/kilo"
$ws.Range("O65").Value = "Región Metropolitana"
$ws.Range("P65").Value = 1000
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"
